$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the member record for "Warqii Lammaa" (row 3) — the whole row
# shifts the remaining members up by one.
$ws.Rows.Item(3).Delete()

# Reflect the post-edit selection location seen in the saved file.
$ws.Range("C12").Select()
